$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new log entry
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A5").Value = "2/9/2020"
$ws.Range("B5").Value = "jaclemon"
$ws.Range("C5").Value = "30 minutes"
$ws.Range("D5").Value = "Used notepad++ to try and diagnose segmentation fault error with csvreader"

# Row 6: new log entry (date+name combined as text in A6, no B6)
$ws.Range("A6").Value = "2/13/2020jaclemon"
$ws.Range("C6").Value = "30 minutes"
$ws.Range("D6").Value = "Used notepad++ to implement -min, -mean functions"

$ws.Range("D6").Select()
